# Append the CSS snippet as new paragraphs after the last paragraph of the
# document body (right before the sectPr), matching the target diff:
#   Body {
#        Front-size: 10px;
#   }
$d = $word.ActiveDocument

$lines = @("Body {", "     Front-size: 10px;", "}")

foreach ($line in $lines) {
    $r = $d.Content
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $r.Collapse(0)
    $r.Text = $line
}
